$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.02
$ws.Range("H2").Value = 3.65
$ws.Range("I2").Value = 3.7
$ws.Range("K2").Value = 4.4
$ws.Range("P2").Value = 2.82
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 1.75
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.8
$ws.Range("V2").Value = 1.37
$ws.Range("W2").Value = 1.98
$ws.Range("X2").Value = 29
$ws.Range("Y2").Value = 23
$ws.Range("AF2").Value = 16
$ws.Range("AJ2").Value = 26
$ws.Range("AK2").Value = 17
$ws.Range("AN2").Value = 8.4
$ws.Range("V3").Value = 1.21
$ws.Range("F4").Value = 5.2
$ws.Range("H4").Value = 1.59
$ws.Range("I4").Value = 1.68
$ws.Range("K4").Value = 4.9
$ws.Range("N4").Value = 5.6
$ws.Range("G5").Value = 2.44
$ws.Range("H5").Value = 4
$ws.Range("W5").Value = 1.69
$ws.Range("AB5").Value = 15.5
$ws.Range("F6").Value = 1.42
$ws.Range("G6").Value = 1.49
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = 1.41
$ws.Range("P6").Value = 1.69
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 2.42
$ws.Range("AJ6").Value = 14.5
$ws.Range("J7").Value = 8
$ws.Range("S7").Value = 1.56
$ws.Range("W7").Value = 6
$ws.Range("F8").Value = 1.3
$ws.Range("G8").Value = 1.37
$ws.Range("H8").Value = 8.800000000000001
$ws.Range("J8").Value = 5.7
$ws.Range("K8").Value = 6.8
$ws.Range("P8").Value = 2.7
$ws.Range("T8").Value = 1.83
$ws.Range("W8").Value = 3.7
$ws.Range("F9").Value = 3.35
$ws.Range("G9").Value = 3.8
$ws.Range("H9").Value = 2.14
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 3.9
$ws.Range("O9").Value = 1.3
$ws.Range("P9").Value = 1.98
$ws.Range("Q9").Value = 1.91
$ws.Range("R9").Value = 1.37
$ws.Range("T9").Value = 1.73
$ws.Range("V9").Value = 1.76
$ws.Range("W9").Value = 1.35
$ws.Range("AJ9").Value = 75
$ws.Range("H10").Value = 2.98
$ws.Range("J10").Value = 3.65
$ws.Range("F11").Value = 2.08
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.65
$ws.Range("L11").Value = 1.19
$ws.Range("R11").Value = 1.69
$ws.Range("S11").Value = 2.2
$ws.Range("T11").Value = 1.5
$ws.Range("W11").Value = 1.83
$ws.Range("X11").Value = 990
$ws.Range("G12").Value = 2.62
$ws.Range("H12").Value = 3.1
$ws.Range("L12").Value = 1.45
$ws.Range("O12").Value = 1.37
$ws.Range("P12").Value = 1.75
$ws.Range("Q12").Value = 1.98
$ws.Range("W12").Value = 1.62
$ws.Range("F13").Value = 7
$ws.Range("H13").Value = 1.43
$ws.Range("I13").Value = 1.49
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 5.9
$ws.Range("L13").Value = 1.2
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 2.52
$ws.Range("S13").Value = 2.34
$ws.Range("T13").Value = 1.77
$ws.Range("V13").Value = 3
$ws.Range("Y13").Value = 12
$ws.Range("AA13").Value = 15.5
$ws.Range("AB13").Value = 32
$ws.Range("F14").Value = 2.2
$ws.Range("G14").Value = 2.38
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 3.5
$ws.Range("J14").Value = 3.65
$ws.Range("O14").Value = 1.24
$ws.Range("T14").Value = 1.65
$ws.Range("V14").Value = 1.4
$ws.Range("AM14").Value = 90
$ws.Range("G15").Value = 1.26
$ws.Range("H15").Value = 12
$ws.Range("I15").Value = 16
$ws.Range("N15").Value = 7.2
$ws.Range("Q15").Value = 1.35
$ws.Range("R15").Value = 1.89
$ws.Range("S15").Value = 1.94
$ws.Range("AJ15").Value = 10.5
$ws.Range("AO15").Value = 220
$ws.Range("G16").Value = 2.3
$ws.Range("AA16").Value = 80
$ws.Range("AE16").Value = 48
$ws.Range("AO16").Value = 42
$ws.Range("G17").Value = 6.6
$ws.Range("H17").Value = 1.64
$ws.Range("I17").Value = 1.7
$ws.Range("K17").Value = 4.7
$ws.Range("Q18").Value = 2.02
$ws.Range("T18").Value = 1.79
$ws.Range("AM18").Value = 120
$ws.Range("T19").Value = 1.46
$ws.Range("F20").Value = 5.9
$ws.Range("I20").Value = 1.57
$ws.Range("K20").Value = 5.9
$ws.Range("V20").Value = 2.78
$ws.Range("AO20").Value = 4.9
$ws.Range("F21").Value = 5.1
$ws.Range("L21").Value = 1.44
$ws.Range("T21").Value = 1.05
$ws.Range("F22").Value = 2.46
$ws.Range("G22").Value = 2.54
$ws.Range("H22").Value = 3.55
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 3.05
$ws.Range("F23").Value = 2.16
$ws.Range("G23").Value = 2.38
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 3.4
$ws.Range("S23").Value = 2.44
$ws.Range("W23").Value = 1.72
$ws.Range("AA23").Value = 65
$ws.Range("AB23").Value = 17.5
$ws.Range("AC23").Value = 10.5
$ws.Range("AF23").Value = 21
$ws.Range("AG23").Value = 14
$ws.Range("AJ23").Value = 36
$ws.Range("AN23").Value = 14
$ws.Range("AO23").Value = 26
$ws.Range("G24").Value = 1.22
$ws.Range("I24").Value = 16
$ws.Range("J24").Value = 8.6
$ws.Range("K24").Value = 10.5
$ws.Range("Q24").Value = 1.27
$ws.Range("S24").Value = 1.68
$ws.Range("T24").Value = 1.73
$ws.Range("U24").Value = 2.1
$ws.Range("W24").Value = 5.4
$ws.Range("X24").Value = 80
$ws.Range("Y24").Value = 100
$ws.Range("AB24").Value = 19
$ws.Range("AF24").Value = 12
$ws.Range("AJ24").Value = 12
$ws.Range("AN24").Value = 2.8
$ws.Range("F25").Value = 2.62
$ws.Range("G25").Value = 2.64
$ws.Range("H25").Value = 2.8
$ws.Range("I25").Value = 2.82
$ws.Range("L25").Value = 1.33
$ws.Range("R25").Value = 1.58
$ws.Range("U25").Value = 2.72
$ws.Range("V25").Value = 1.54
$ws.Range("W25").Value = 1.6
$ws.Range("AB25").Value = 14
$ws.Range("AK25").Value = 24
$ws.Range("AN25").Value = 17
$ws.Range("AO25").Value = 18.5
